# Automatic update: dades i banners [2026-02-21 00:19]
# Advances the daily meteocat observation date from 2026-02-20 to 2026-02-21,
# refreshes the per-row extraction timestamp / source URL, and replaces the
# (now stale) observed values with "sense dades" placeholders since the new
# day has no data yet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldDate = "2026-02-20"
$newDate = "2026-02-21"

# New DATA_EXTRACCIO timestamp per data row (row number -> "yyyy-MM-dd HH:mm:ss")
$newExtractionTime = @{
    2  = "2026-02-21 00:18:18"
    3  = "2026-02-21 00:18:20"
    4  = "2026-02-21 00:18:22"
    5  = "2026-02-21 00:18:24"
    6  = "2026-02-21 00:18:26"
    7  = "2026-02-21 00:18:28"
    8  = "2026-02-21 00:18:31"
    9  = "2026-02-21 00:18:32"
    10 = "2026-02-21 00:18:35"
    11 = "2026-02-21 00:18:36"
    12 = "2026-02-21 00:18:37"
    13 = "2026-02-21 00:18:38"
    14 = "2026-02-21 00:18:39"
    15 = "2026-02-21 00:18:40"
    16 = "2026-02-21 00:18:41"
    17 = "2026-02-21 00:18:42"
    18 = "2026-02-21 00:18:43"
    19 = "2026-02-21 00:18:45"
    20 = "2026-02-21 00:18:46"
    21 = "2026-02-21 00:18:47"
    22 = "2026-02-21 00:18:48"
    23 = "2026-02-21 00:18:49"
    24 = "2026-02-21 00:18:51"
    25 = "2026-02-21 00:18:53"
    26 = "2026-02-21 00:18:56"
    27 = "2026-02-21 00:18:58"
    28 = "2026-02-21 00:19:00"
    29 = "2026-02-21 00:19:02"
    30 = "2026-02-21 00:19:04"
    31 = "2026-02-21 00:19:06"
    32 = "2026-02-21 00:19:09"
    33 = "2026-02-21 00:19:11"
    34 = "2026-02-21 00:19:13"
    35 = "2026-02-21 00:19:15"
    36 = "2026-02-21 00:19:17"
    37 = "2026-02-21 00:19:19"
    38 = "2026-02-21 00:19:21"
    39 = "2026-02-21 00:19:23"
    40 = "2026-02-21 00:19:25"
    41 = "2026-02-21 00:19:27"
    42 = "2026-02-21 00:19:29"
    43 = "2026-02-21 00:19:31"
    44 = "2026-02-21 00:19:33"
    45 = "2026-02-21 00:19:35"
    46 = "2026-02-21 00:19:37"
}

# Columns (1-indexed): A=1 ID_ESTAC, D=4 DATA_DIA, E=5 DATA_EXTRACCIO, F=6 URL_FONT,
# G=7 GRUIX_NEU_MAX, H=8 HUMITAT_MITJANA_DIA, I=9 PRECIPITACIO_ACUM_DIA,
# J=10 PRESSIO_ATMOSFERICA, K=11 RADIACIO_GLOBAL, L=12 RATXA_VENT_MAX,
# M=13 TEMPERATURA_MAXIMA_DIA, N=14 TEMPERATURA_MINIMA_DIA, O=15 TEMPERATURA_MITJANA_DIA

$firstDataRow = 2
$lastDataRow = 46

# DATA_DIA (column D) holds a plain yyyy-MM-dd string; force text storage up
# front so Excel doesn't reinterpret the assigned string as a real date.
$ws.Range("D$firstDataRow`:D$lastDataRow").NumberFormat = "@"

for ($r = $firstDataRow; $r -le $lastDataRow; $r++) {

    $codi = $ws.Cells.Item($r, 1).Value2

    # D: DATA_DIA
    $ws.Cells.Item($r, 4).Value = $newDate

    # E: DATA_EXTRACCIO
    $ws.Cells.Item($r, 5).Value = $newExtractionTime[$r]

    # F: URL_FONT - same station code, date bumped by one day
    $ws.Cells.Item($r, 6).Value = "https://www.meteo.cat/observacions/xema/dades?codi=$codi&dia=$newDate" + "T09:00Z"

    # G: GRUIX_NEU_MAX -> "sense dades" if it previously had a value
    $gVal = $ws.Cells.Item($r, 7).Value2
    if (-not [string]::IsNullOrEmpty($gVal)) {
        $ws.Cells.Item($r, 7).Value = "sense dades"
    }

    # H: HUMITAT_MITJANA_DIA -> "sense dades" if it previously had a value
    $hVal = $ws.Cells.Item($r, 8).Value2
    if (-not [string]::IsNullOrEmpty($hVal)) {
        $ws.Cells.Item($r, 8).Value = "sense dades"
    }

    # I: PRECIPITACIO_ACUM_DIA -> "sense dades" if it previously had a value
    $iVal = $ws.Cells.Item($r, 9).Value2
    if (-not [string]::IsNullOrEmpty($iVal)) {
        $ws.Cells.Item($r, 9).Value = "sense dades"
    }

    # J: PRESSIO_ATMOSFERICA -> always cleared (no longer reported)
    $jVal = $ws.Cells.Item($r, 10).Value2
    if (-not [string]::IsNullOrEmpty($jVal)) {
        $ws.Cells.Item($r, 10).Value = ""
    }

    # K: RADIACIO_GLOBAL -> "sense dades" if it previously had a value
    $kVal = $ws.Cells.Item($r, 11).Value2
    if (-not [string]::IsNullOrEmpty($kVal)) {
        $ws.Cells.Item($r, 11).Value = "sense dades"
    }

    # L: RATXA_VENT_MAX -> "sense dades sense dades" if it previously had a value
    $lVal = $ws.Cells.Item($r, 12).Value2
    if (-not [string]::IsNullOrEmpty($lVal)) {
        $ws.Cells.Item($r, 12).Value = "sense dades sense dades"
    }

    # M: TEMPERATURA_MAXIMA_DIA -> "sense dades sense dades" if it previously had a value
    $mVal = $ws.Cells.Item($r, 13).Value2
    if (-not [string]::IsNullOrEmpty($mVal)) {
        $ws.Cells.Item($r, 13).Value = "sense dades sense dades"
    }

    # N: TEMPERATURA_MINIMA_DIA -> "sense dades sense dades" if it previously had a value
    $nVal = $ws.Cells.Item($r, 14).Value2
    if (-not [string]::IsNullOrEmpty($nVal)) {
        $ws.Cells.Item($r, 14).Value = "sense dades sense dades"
    }

    # O: TEMPERATURA_MITJANA_DIA -> "sense dades" if it previously had a value
    $oVal = $ws.Cells.Item($r, 15).Value2
    if (-not [string]::IsNullOrEmpty($oVal)) {
        $ws.Cells.Item($r, 15).Value = "sense dades"
    }
}

Write-Host "Update complete: rows $firstDataRow to $lastDataRow advanced from $oldDate to $newDate"
